$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Comments text for the BMP280 row (row 3)
$ws.Range("D3").Value = "Pressure & Altitude  & Temperature"

# Add a new "Comments" column G, copying formatting from column F
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("F2").Copy($ws.Range("G2:G7"))

$ws.Range("G1").Value = "Comments"
$ws.Range("G3").Value = "Altitude, Temperature not sending."

# Size the new column
$ws.Columns.Item(7).ColumnWidth = 16.42578125

# Update the active selection
$ws.Range("G4").Select()
